$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 777
$ws.Range("F3").Value = 2815
$ws.Range("F4").Value = 1340
$ws.Range("F6").Value = 590
$ws.Range("F7").Value = 49
$ws.Range("F9").Value = 284
$ws.Range("F11").Value = 11734
$ws.Range("F12").Value = 6667
$ws.Range("F14").Value = 17
$ws.Range("F15").Value = 423
$ws.Range("F16").Value = 255
$ws.Range("F18").Value = 12
$ws.Range("F19").Value = 926
$ws.Range("F20").Value = 91
$ws.Range("F21").Value = 278
$ws.Range("F22").Value = 931
$ws.Range("F23").Value = 3653
$ws.Range("F25").Value = 988
$ws.Range("F26").Value = 498
$ws.Range("F27").Value = 173
$ws.Range("F30").Value = 229
$ws.Range("F32").Value = 309
$ws.Range("F33").Value = 5030
$ws.Range("F34").Value = 44
$ws.Range("F35").Value = 1246
$ws.Range("F36").Value = 239
$ws.Range("F37").Value = 525
$ws.Range("F38").Value = 206

$ws = $wb.Worksheets.Item(2)
$ws.Range("F5").Value = 8
$ws.Range("F11").Value = 3685

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 9064
$ws.Range("F3").Value = 507
$ws.Range("F4").Value = 1835

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 507
$ws.Range("F3").Value = 1835
$ws.Range("F4").Value = 777
$ws.Range("F5").Value = 2815
$ws.Range("F9").Value = 1340
$ws.Range("F11").Value = 590
$ws.Range("F12").Value = 49
$ws.Range("F15").Value = 284
$ws.Range("F17").Value = 11734
$ws.Range("F18").Value = 3685
$ws.Range("F19").Value = 6667
$ws.Range("F22").Value = 17
$ws.Range("F23").Value = 423
$ws.Range("F24").Value = 255
$ws.Range("F26").Value = 12
$ws.Range("F27").Value = 91
$ws.Range("F28").Value = 278
$ws.Range("F29").Value = 931
$ws.Range("F30").Value = 3653
$ws.Range("F32").Value = 988
$ws.Range("F33").Value = 173
$ws.Range("F35").Value = 229
$ws.Range("F39").Value = 44
$ws.Range("F40").Value = 1246
$ws.Range("F41").Value = 239
$ws.Range("F43").Value = 206
